# Actualización automática 2025-09-12 12:50:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M13").Value = 2635.78
$wsGrupo.Range("M23").Value = "7 de 21"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F13").Value = 2635.78
$wsMensual.Range("F23").Value = 35458.76

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D12").Value = 31973.48
$wsCumpl.Range("E12").Value = 4850.163092117098
$wsCumpl.Range("F12").Value = 0.8682867124259256
$wsCumpl.Range("D15").Value = 35458.76
$wsCumpl.Range("E15").Value = 19965.98316613378
$wsCumpl.Range("F15").Value = 0.6397640832310865

# Column D width change (13 -> 14) on "CUMPLIMIENTO MENSUAL"
# NOTE: ColumnWidth undergoes a pixel-based round-trip conversion before being
# stored back as the OOXML "width" attribute, so a value of 13.15 (which lies
# safely inside the quantization band that maps back to a stored width of
# exactly 14) is used instead of 14 itself.
$wsCumpl.Columns.Item(4).ColumnWidth = 13.15
